$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.415.98'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.893.08'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4892'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2930'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06706'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.905.77'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.98'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07343'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.126'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '87.62'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6635'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.381.95'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.42'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007828'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.326'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +12.15%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.107.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.33%  '
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '189.65'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.112'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.467'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.929'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.460'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.345'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09142'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.039'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05188'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7389'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.099'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.716'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.02%  '
$ws.Range('E37').Value = '  -0.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.663'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9209'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.038'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.32%  '
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.920'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '105.98'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9924'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.71%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1371'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.73%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '68.20'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +19.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.586'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.958'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.85'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05818'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3938'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.46%  '
